# Update "想去人数" (interested-count) figures across sheets to reflect
# a refreshed data scrape (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 54349
$ws1.Range("F7").Value = 1316
$ws1.Range("F8").Value = 335
$ws1.Range("F9").Value = 311
$ws1.Range("F11").Value = 722
$ws1.Range("F13").Value = 2993
$ws1.Range("F14").Value = 876
$ws1.Range("F16").Value = 1260
$ws1.Range("F19").Value = 263
$ws1.Range("F21").Value = 381
$ws1.Range("F22").Value = 1234
$ws1.Range("F23").Value = 87
$ws1.Range("F25").Value = 162
$ws1.Range("F30").Value = 63
$ws1.Range("F32").Value = 4798
$ws1.Range("F34").Value = 4805
$ws1.Range("F35").Value = 8773
$ws1.Range("F38").Value = 127
$ws1.Range("F42").Value = 74
$ws1.Range("F44").Value = 199

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 92

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 13212
$ws4.Range("F7").Value = 1316
$ws4.Range("F8").Value = 335
$ws4.Range("F9").Value = 311
$ws4.Range("F11").Value = 722
$ws4.Range("F13").Value = 2993
$ws4.Range("F14").Value = 876
$ws4.Range("F15").Value = 92
$ws4.Range("F16").Value = 1260
$ws4.Range("F22").Value = 263
$ws4.Range("F26").Value = 87
$ws4.Range("F27").Value = 162
$ws4.Range("F31").Value = 4798
$ws4.Range("F33").Value = 4805
$ws4.Range("F37").Value = 127
$ws4.Range("F43").Value = 74

$wb.Save()
